$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D), Volume(1h) (E), and Hora (G) columns for the crypto
# symbol list refresh. Values are stored as plain text in the sheet
# (matching the original inline-string cells), so we force the
# NumberFormat to Text ("@") before assigning each value to avoid
# Excel silently converting numeric-looking / percentage-looking
# strings into Number cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '327.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.87%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '16'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.14%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '16'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.672'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '10.62%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '16'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08122'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.64%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '16'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.579'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.19%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '16'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.717'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.07%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '16'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.953'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.15%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '16'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.95%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '16'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9461'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.63%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '16'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1287'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '13.46%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '16'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1991'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '6.37%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '16'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09245'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.81%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '16'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '6.11%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '16'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09621'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.76%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '16'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001314'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-4.57%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '16'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006299'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3.01%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '16'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.370'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.62%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '16'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3537'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.40%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '16'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.624'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '19.93%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '16'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1419'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '9.79%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '16'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.55%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '16'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04452'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.33%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '16'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001252'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '4.16%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '16'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004337'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.52%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '16'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001191'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.96%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '16'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003989'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '37.25%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '16'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '16'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '16'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '16'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '16'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '16'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '16'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '16'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '16'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '16'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '16'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '16'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02511'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '17.77%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '16'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05254'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '6.50%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '16'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007382'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.80%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '16'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1433'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '5.98%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '16'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008883'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.00%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '16'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002069'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.17%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '16'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01090'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '26.67%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '16'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006765'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.64%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '16'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.14%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '16'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002875'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-12.90%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '16'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001800'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '24.45%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '16'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.14%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '16'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.14%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '16'
